# 🧠 add SpriteSetMesh() logic
# Highlight key bullet points on the "Step 4" slide (Sprite.c) so the
# Create/Free/MeshRender related lines stand out.
#
# Color reference (COM Font.Highlight takes a VBA-style RGB long:
#   value = R + G*256 + B*65536):
#   Green  00FF00 -> 65280
#   Yellow FFFF00 -> 65535

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$green = 65280
$yellow = 65535

# "Implement the Create function"
$tr.Paragraphs(2, 1).Font.Highlight = $green
# "Test to make sure that the object is constructed properly"
$tr.Paragraphs(3, 1).Font.Highlight = $yellow
# "Implement the Free function"
$tr.Paragraphs(4, 1).Font.Highlight = $green
# "Test to make sure that the object is freed properly"
$tr.Paragraphs(5, 1).Font.Highlight = $yellow
# "Test to make sure that the original pointer is set to NULL"
$tr.Paragraphs(6, 1).Font.Highlight = $yellow

# "Use the MeshRender function for this purpose" (paragraph 9)
$full = $tr.Text
$target = "Use the MeshRender function for this purpose"
$idx = $full.IndexOf($target)
$tr.Characters($idx + 1, [int]$target.Length).Font.Highlight = $yellow
